$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B width (14.7109375 -> 15.42578125 characters)
$ws.Columns.Item(2).ColumnWidth = 15.42578125

# Update cell values in columns A and B for rows 1-32
$ws.Cells.Item(1, 1).Value = -0.30448347740663451
$ws.Cells.Item(1, 2).Value = 0.30357134761589322
$ws.Cells.Item(2, 1).Value = -0.23434420572618908
$ws.Cells.Item(2, 2).Value = 0.23158670674378623
$ws.Cells.Item(3, 1).Value = -0.090808514708390575
$ws.Cells.Item(3, 2).Value = 0.090487922041708302
$ws.Cells.Item(4, 1).Value = -0.082487922090740184
$ws.Cells.Item(4, 2).Value = 0.082210768690197611
$ws.Cells.Item(5, 1).Value = -0.079210768718958491
$ws.Cells.Item(5, 2).Value = 0.078287814848184922
$ws.Cells.Item(6, 1).Value = -0.01603403751328969
$ws.Cells.Item(6, 2).Value = 0.015924314797244321
$ws.Cells.Item(7, 1).Value = -0.0059243148620966579
$ws.Cells.Item(7, 2).Value = 0.0059108400774983139
$ws.Cells.Item(8, 1).Value = 0.0040891598574335219
$ws.Cells.Item(8, 2).Value = -0.0040937251791057783
$ws.Cells.Item(9, 1).Value = 0.0060937251506616441
$ws.Cells.Item(9, 2).Value = -0.0060958277885636569
$ws.Cells.Item(10, 1).Value = 0.0080958277608704776
$ws.Cells.Item(10, 2).Value = -0.0080955717409327832
$ws.Cells.Item(11, 1).Value = -0.024393665969102862
$ws.Cells.Item(11, 2).Value = 0.024367135021029718
$ws.Cells.Item(12, 1).Value = -0.020867135056049868
$ws.Cells.Item(12, 2).Value = 0.02067183770096559
$ws.Cells.Item(13, 1).Value = -0.017171837738131579
$ws.Cells.Item(13, 2).Value = 0.01708284834378393
$ws.Cells.Item(14, 1).Value = -0.0090828484017695388
$ws.Cells.Item(14, 2).Value = 0.0090537350118049531
$ws.Cells.Item(15, 1).Value = -0.0080537350386071793
$ws.Cells.Item(15, 2).Value = 0.008034932845781384
$ws.Cells.Item(16, 1).Value = -0.0060349328774682576
$ws.Cells.Item(16, 2).Value = 0.006003553143477447
$ws.Cells.Item(17, 1).Value = -0.0040035531756279497
$ws.Cells.Item(17, 2).Value = 0.0039999999586770585
$ws.Cells.Item(18, 1).Value = -0.063227141918481777
$ws.Cells.Item(18, 2).Value = 0.063069837458098021
$ws.Cells.Item(19, 1).Value = -0.012092158704574363
$ws.Cells.Item(19, 2).Value = 0.012016856413185728
$ws.Cells.Item(20, 1).Value = -0.008016856435041575
$ws.Cells.Item(20, 2).Value = 0.0080056861172206339
$ws.Cells.Item(21, 1).Value = -0.0040056861392976373
$ws.Cells.Item(21, 2).Value = 0.0039999999777498019
$ws.Cells.Item(22, 1).Value = -0.136778192064134
$ws.Cells.Item(22, 2).Value = 0.1356384846541383
$ws.Cells.Item(23, 1).Value = -0.10120342106467017
$ws.Cells.Item(23, 2).Value = 0.099909281553330231
$ws.Cells.Item(24, 1).Value = -0.079909281666067145
$ws.Cells.Item(24, 2).Value = 0.079235648792076674
$ws.Cells.Item(25, 1).Value = -0.060253777363817207
$ws.Cells.Item(25, 2).Value = 0.060204809846368335
$ws.Cells.Item(26, 1).Value = -0.057704809878300622
$ws.Cells.Item(26, 2).Value = 0.05764464685750248
$ws.Cells.Item(27, 1).Value = -0.055144646890664895
$ws.Cells.Item(27, 2).Value = 0.054805933580905908
$ws.Cells.Item(28, 1).Value = -0.052805933616621559
$ws.Cells.Item(28, 2).Value = 0.05258812996636486
$ws.Cells.Item(29, 1).Value = -0.045588130027940821
$ws.Cells.Item(29, 2).Value = 0.045535809589376086
$ws.Cells.Item(30, 1).Value = 0.014464190108755659
$ws.Cells.Item(30, 2).Value = -0.014484168913614592
$ws.Cells.Item(31, 1).Value = -0.014023484540681252
$ws.Cells.Item(31, 2).Value = 0.01400130359244578
$ws.Cells.Item(32, 1).Value = -0.0040013036683728131
$ws.Cells.Item(32, 2).Value = 0.003999999950998756
